# Auto-generated Excel COM-interop script
# Updates market-price / profit columns (H:N) for specific leve rows
# across multiple crafting-job sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3828.697
$ws.Range("J17").Value = 3828.697
$ws.Range("L17").Value = 11486.091
$ws.Range("N17").Value = -11822.091

$ws.Range("H112").Value = 2492.125
$ws.Range("I112").Value = 814.6667
$ws.Range("J112").Value = 2879.2307
$ws.Range("K112").Value = 2444.0001
$ws.Range("L112").Value = 8637.6921
$ws.Range("M112").Value = -1336.0001
$ws.Range("N112").Value = -10853.6921

$ws.Range("H129").Value = 1045.8586
$ws.Range("J129").Value = 973.3617
$ws.Range("L129").Value = 2920.0851
$ws.Range("N129").Value = -12920.0851

$ws.Range("H138").Value = 1254.56
$ws.Range("I138").Value = 572.93616
$ws.Range("J138").Value = 1859.0189
$ws.Range("K138").Value = 1718.80848
$ws.Range("L138").Value = 5577.0567
$ws.Range("M138").Value = 3421.19152
$ws.Range("N138").Value = -15857.0567

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3378.4443
$ws.Range("I61").Value = 2005.75
$ws.Range("K61").Value = 2005.75
$ws.Range("M61").Value = -1793.75

$ws.Range("H102").Value = 73969.836
$ws.Range("I102").Value = 3799
$ws.Range("J102").Value = 88004
$ws.Range("K102").Value = 3799
$ws.Range("L102").Value = 88004
$ws.Range("M102").Value = -2177
$ws.Range("N102").Value = -91248

$ws.Range("H136").Value = 3378.4443
$ws.Range("I136").Value = 2005.75
$ws.Range("K136").Value = 6017.25
$ws.Range("M136").Value = -3467.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2467.3333
$ws.Range("I105").Value = 1803
$ws.Range("J105").Value = 3614.818
$ws.Range("K105").Value = 1803
$ws.Range("L105").Value = 3614.818
$ws.Range("M105").Value = -56
$ws.Range("N105").Value = -7108.818

$ws.Range("H107").Value = 1800.8125
$ws.Range("I107").Value = 1420.45
$ws.Range("J107").Value = 2434.75
$ws.Range("K107").Value = 1420.45
$ws.Range("L107").Value = 2434.75
$ws.Range("M107").Value = 499.55
$ws.Range("N107").Value = -6274.75

$ws.Range("H109").Value = 19998.4
$ws.Range("J109").Value = 19998.4
$ws.Range("L109").Value = 19998.4
$ws.Range("N109").Value = -22772.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4844.6577
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4844.6577
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4844.6577
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -5434.6577

$ws.Range("H34").Value = 4844.6577
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4844.6577
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 4844.6577
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5248.6577

$ws.Range("H58").Value = 1517.0256
$ws.Range("I58").Value = 1225.9412
$ws.Range("J58").Value = 3496.4
$ws.Range("K58").Value = 1225.9412
$ws.Range("L58").Value = 3496.4
$ws.Range("M58").Value = -1022.9412
$ws.Range("N58").Value = -3902.4

$ws.Range("H99").Value = 3135.1667
$ws.Range("J99").Value = 3324.75
$ws.Range("L99").Value = 3324.75
$ws.Range("N99").Value = -6320.75

$ws.Range("H126").Value = 3135.1667
$ws.Range("J126").Value = 3324.75
$ws.Range("L126").Value = 9974.25
$ws.Range("N126").Value = -14914.25

$ws.Range("H136").Value = 1517.0256
$ws.Range("I136").Value = 1225.9412
$ws.Range("J136").Value = 3496.4
$ws.Range("K136").Value = 3677.8236
$ws.Range("L136").Value = 10489.2
$ws.Range("M136").Value = -1127.8236
$ws.Range("N136").Value = -15589.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3431.2974
$ws.Range("I5").Value = 4302.615
$ws.Range("K5").Value = 12907.845
$ws.Range("M5").Value = -12795.845

$ws.Range("H36").Value = 1469.75
$ws.Range("I36").Value = 293.33334
$ws.Range("J36").Value = 4999
$ws.Range("K36").Value = 880.0000200000001
$ws.Range("L36").Value = 14997
$ws.Range("M36").Value = -711.0000200000001
$ws.Range("N36").Value = -15335

$ws.Range("H70").Value = 4027.7273
$ws.Range("I70").Value = 2087.7144
$ws.Range("J70").Value = 4933.067
$ws.Range("K70").Value = 6263.1432
$ws.Range("L70").Value = 14799.201
$ws.Range("M70").Value = -5948.1432
$ws.Range("N70").Value = -15429.201

$ws.Range("H73").Value = 4027.7273
$ws.Range("I73").Value = 2087.7144
$ws.Range("J73").Value = 4933.067
$ws.Range("K73").Value = 6263.1432
$ws.Range("L73").Value = 14799.201
$ws.Range("M73").Value = -5171.1432
$ws.Range("N73").Value = -16983.201

$ws.Range("H104").Value = 1500
$ws.Range("I104").Value = 1500
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 4500
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -1879
$ws.Range("N104").ClearContents()

$ws.Range("H113").Value = 6000.7896
$ws.Range("I113").Value = 10434.1
$ws.Range("J113").Value = 1074.8889
$ws.Range("K113").Value = 31302.3
$ws.Range("L113").Value = 3224.6667
$ws.Range("M113").Value = -29132.3
$ws.Range("N113").Value = -7564.6667

$ws.Range("H121").Value = 130096.625
$ws.Range("I121").Value = 615
$ws.Range("J121").Value = 141867.69
$ws.Range("K121").Value = 1845
$ws.Range("L121").Value = 425603.07
$ws.Range("M121").Value = -535
$ws.Range("N121").Value = -428223.07

$ws.Range("H135").Value = 3431.2974
$ws.Range("I135").Value = 4302.615
$ws.Range("K135").Value = 38723.535
$ws.Range("M135").Value = -36188.535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2700.3333
$ws.Range("I102").Value = 2752.875
$ws.Range("K102").Value = 2752.875
$ws.Range("M102").Value = -1130.875

$ws.Range("H113").Value = 1506.75
$ws.Range("I113").Value = 1414.8572
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 1414.8572
$ws.Range("L113").Value = 2150
$ws.Range("M113").Value = 755.1428000000001
$ws.Range("N113").Value = -6490

$ws.Range("H122").Value = 1563.8889
$ws.Range("I122").Value = 1018.75
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3056.25
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -606.25
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3612.7144
$ws.Range("I40").Value = 2236.6365
$ws.Range("J40").Value = 8658.333000000001
$ws.Range("K40").Value = 2236.6365
$ws.Range("L40").Value = 8658.333000000001
$ws.Range("M40").Value = -2100.6365
$ws.Range("N40").Value = -8930.333000000001

$ws.Range("H93").Value = 1339.5
$ws.Range("I93").Value = 881.7692
$ws.Range("K93").Value = 881.7692
$ws.Range("M93").Value = 366.2308

$ws.Range("H100").Value = 2072.9167
$ws.Range("I100").Value = 2006.8182
$ws.Range("K100").Value = 2006.8182
$ws.Range("M100").Value = -1465.8182

$ws.Range("H132").Value = 2601.7036
$ws.Range("I132").Value = 2031.6578
$ws.Range("J132").Value = 3955.5625
$ws.Range("K132").Value = 6094.9734
$ws.Range("L132").Value = 11866.6875
$ws.Range("M132").Value = -3564.9734
$ws.Range("N132").Value = -16926.6875
